$d = $word.ActiveDocument
$pBreak = $d.Paragraphs.Item(91)
$pBreak.Range.Delete()
$p = $d.Paragraphs.Item(91)
Write-Output ("p91 text=[" + $p.Range.Text + "]")

# Now re-type "C" at the start to force a text edit / run regen.
$rC = $d.Range($p.Range.Start, $p.Range.Start + 1)
$rC.Text = "C"
$p2 = $d.Paragraphs.Item(91)
Write-Output ("p91 now text=[" + $p2.Range.Text + "]")

# Now try to re-split the run: select "AREFUL NOTE:" substring (length 12) and
# re-apply formatting distinctly.
$rest = $d.Range($p2.Range.Start + 1, $p2.Range.Start + 1 + 12)
Write-Output ("rest text=[" + $rest.Text + "]")
$rest.Font.Color = 255  # blue, to force distinguishable formatting
Write-Output "set rest color blue"
$p3 = $d.Paragraphs.Item(91)
Write-Output ("p91 now text=[" + $p3.Range.Text + "]")

# Now set it back to red to match original, see if runs stay split or merge again.
$rest2 = $d.Range($p3.Range.Start + 1, $p3.Range.Start + 1 + 12)
$rest2.Font.Color = 255
Write-Output ("rest2 color reset attempt")
